$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# --- Shape 4: "Direct Measurement TLV" content placeholder ---
$shape1 = $s.Shapes.Item(4)
$tr1 = $shape1.TextFrame.TextRange

# Add the new 4th paragraph
$tr1.InsertAfter("`rDoes not support per-traffic class direct measurement") | Out-Null

# --- Shape 5: "Direct Measurement Test Packet" content placeholder ---
$shape2 = $s.Shapes.Item(5)
$tr2 = $shape2.TextFrame.TextRange

# Update "packet and byte " -> "packet and 64-bit byte " (4th run of 3rd paragraph)
$para3 = $tr2.Paragraphs(3)
$run4 = $para3.Runs(4)
$run4.Text = "packet and 64-bit byte "

# Add the new 5th paragraph
$tr2.InsertAfter("`rPlan to add: per traffic-class counter collection (per traffic-class loss measurement) (Ok to drop best effort traffic)") | Out-Null

# --- Resize all text in both boxes from 16pt to 14pt ---
foreach ($shape in @($shape1, $shape2)) {
    $tr = $shape.TextFrame.TextRange
    $nParas = $tr.Paragraphs().Count
    for ($i = 1; $i -le $nParas; $i++) {
        $para = $tr.Paragraphs($i)
        $nRuns = $para.Runs().Count
        for ($j = 1; $j -le $nRuns; $j++) {
            $run = $para.Runs($j)
            $run.Font.Size = 14
        }
    }
}
